$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.336.20'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  +1.61%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.325.26'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("E4").Value = '  +0.09%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.10'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +5.80%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.43'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +1.67%  '

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.993'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.54%  '

$ws.Range("E8").Value = '  +0.55%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.352.42'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +0.69%  '

$ws.Range("E10").Value = '  +1.98%  '

$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("E12").Value = '  +1.26%  '

$ws.Range("E13").Value = '  +4.70%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.750.87'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +0.17%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.50'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -0.16%  '

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.422.04'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +1.77%  '

$ws.Range("E17").Value = '  +0.70%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.352.43'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +0.90%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '336.18'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +3.69%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.54'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("E21").Value = '  +1.57%  '

$ws.Range("E22").Value = '  +1.90%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.00%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.78'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +0.15%  '

$ws.Range("E25").Value = '  +4.32%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -0.17%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.44'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -2.64%  '

$ws.Range("E28").Value = '  +8.09%  '

$ws.Range("E29").Value = '  +4.68%  '

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.59'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +1.81%  '

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0734'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +1.99%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.16'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.98%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.56'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +1.48%  '

$ws.Range("E34").Value = '  +14.41%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -0.10%  '

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E38").Value = '  +4.25%  '

$ws.Range("E39").Value = '  +3.86%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.37'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +2.40%  '

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.51'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("E43").Value = '  +1.21%  '

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '281.26'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +0.67%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.27'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +6.78%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0930'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +0.57%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0505'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +2.02%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.560'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("E49").Value = '  +1.91%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.382'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("E51").Value = '  +1.68%  '
